$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).ClearFormats()
}

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}


# Row 2 - Bitcoin
Set-Cell "D2" "66.911.79"
Set-Cell "E2" "  -0.48%  "

# Row 3 - Ethereum
Set-Cell "D3" "3.114.24"
Set-Cell "E3" "  +0.27%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  +0.07%  "

# Row 5 - BNB
Set-TextCell "D5" "577.52"
Set-Cell "E5" "  -0.94%  "

# Row 6 - Solana
Set-TextCell "D6" "171.58"
Set-Cell "E6" "  +0.66%  "

# Row 7 - USDC
Set-Cell "E7" "  +0.07%  "

# Row 8 - LidoStakedEther
Set-Cell "D8" "3.110.82"
Set-Cell "E8" "  +0.29%  "

# Row 9 - XRP
Set-Cell "E9" "  -1.03%  "

# Row 10 - Toncoin
Set-Cell "E10" "  -3.25%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.152"
Set-Cell "E11" "  -1.97%  "

# Row 12 - Cardano
Set-TextCell "D12" "0.482"
Set-Cell "E12" "  -0.15%  "

# Row 13 - ShibaInu
Set-TextCell "D13" "0.0000245"
Set-Cell "E13" "  -2.42%  "

# Row 14 - Avalanche
Set-TextCell "D14" "37.25"
Set-Cell "E14" "  +0.56%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-Cell "D16" "3.633.76"
Set-Cell "E16" "  +0.39%  "

# Row 17 - WrappedBTC
Set-Cell "D17" "66.867.38"
Set-Cell "E17" "  -0.48%  "

# Row 18 - Polkadot
Set-TextCell "D18" "7.13"
Set-Cell "E18" "  -1.81%  "

# Row 19 - WrappedEther
Set-Cell "D19" "3.116.65"
Set-Cell "E19" "  +0.38%  "

# Row 20 - Chainlink
Set-TextCell "D20" "16.44"
Set-Cell "E20" "  +0.86%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "475.77"
Set-Cell "E21" "  +0.69%  "

# Row 22 - was Uniswap, becomes Polygon
Set-Cell "B22" "Polygon"
Set-Cell "C22" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D22" "0.713"
Set-Cell "E22" "  -0.63%  "

# Row 23 - was Polygon, becomes Uniswap
Set-Cell "B23" "Uniswap"
Set-Cell "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D23" "7.95"
Set-Cell "E23" "  +5.04%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextCell "D24" "13.49"
Set-Cell "E24" "  +4.19%  "

# Row 25 - Litecoin
Set-TextCell "D25" "83.98"
Set-Cell "E25" "  -0.09%  "

# Row 26 - Fetch.AI
Set-TextCell "D26" "2.30"
Set-Cell "E26" "  -3.06%  "

# Row 27 - Dai
Set-TextCell "D27" "0.999"
Set-Cell "E27" "  -0.03%  "

# Row 28 - RenderToken
Set-TextCell "D28" "10.00"
Set-Cell "E28" "  -2.48%  "

# Row 29 - NEARProtocol
Set-TextCell "D29" "7.88"
Set-Cell "E29" "  -2.90%  "

# Row 30 - ImmutableX
Set-Cell "E30" "  -1.62%  "

# Row 31 - PancakeSwap
Set-Cell "E31" "  -0.63%  "

# Row 32 - EthereumClassic
Set-TextCell "D32" "28.52"
Set-Cell "E32" "  +0.08%  "

# Row 33 - Hedera
Set-Cell "E33" "  -0.43%  "

# Row 34 - PEPE
Set-Cell "D34" "0.0₃0935"
Set-Cell "E34" "  -8.48%  "

# Row 35 - FirstDigitalUSD
Set-TextCell "D35" "1.00"
Set-Cell "E35" "  +0.10%  "

# Row 36 - Filecoin
Set-TextCell "D36" "5.86"
Set-Cell "E36" "  -1.46%  "

# Row 37 - Mantle
Set-TextCell "D37" "0.978"
Set-Cell "E37" "  -3.77%  "

# Row 38 - Arweave
Set-TextCell "D38" "47.26"
Set-Cell "E38" "  -0.28%  "

# Row 39 - Stacks
Set-Cell "E39" "  -1.43%  "

# Row 40 - OKB
Set-TextCell "D40" "50.01"
Set-Cell "E40" "  -0.81%  "

# Row 41 - TheGraph
Set-Cell "E41" "  -2.64%  "

# Row 42 - Kaspa
Set-Cell "E42" "  -1.68%  "

# Row 43 - Cosmos
Set-TextCell "D43" "8.68"
Set-Cell "E43" "  -0.89%  "

# Row 44 - Maker
Set-Cell "D44" "2.807.28"
Set-Cell "E44" "  +1.23%  "

# Row 45 - VeChain
Set-TextCell "D45" "0.0356"
Set-Cell "E45" "  -2.39%  "

# Row 46 - Bittensor
Set-TextCell "D46" "381.54"
Set-Cell "E46" "  -3.87%  "

# Row 47 - dogwifhat
Set-Cell "E47" "  -12.17%  "

# Row 48 - Monero
Set-TextCell "D48" "136.05"
Set-Cell "E48" "  +0.53%  "

# Row 49 - USDe
Set-Cell "E49" "  +0.04%  "

# Row 50 - InjectiveProtocol
Set-TextCell "D50" "24.67"
Set-Cell "E50" "  -0.78%  "

# Row 51 - ThetaToken
Set-Cell "E51" "  -2.41%  "
